$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Shahbaz Nadeem"

# Insert a new column before column A, shifting all existing data right by one
$ws.Columns("A").Insert()

# Populate the new column A with the "matchNo" header/value
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "6th"
